$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values of B7 and B8
$b7 = $ws.Range("B7").Value2
$b8 = $ws.Range("B8").Value2
$ws.Range("B7").Value = $b8
$ws.Range("B8").Value = $b7

# Swap values of B29 and B30
$b29 = $ws.Range("B29").Value2
$b30 = $ws.Range("B30").Value2
$ws.Range("B29").Value = $b30
$ws.Range("B30").Value = $b29

# Give B29 an explicit black font color (RGB 0,0,0)
$ws.Range("B29").Font.Color = 0

# Add a formatted (empty) cell at C30, copying B30's formatting
$ws.Range("B30").Copy()
$ws.Range("C30").PasteSpecial(-4122)

# Update the active selection to C30
$ws.Range("C30").Select()
